# chore: update Sheets via scheduled runner
# Refresh computed Leve-profit columns (currentAveragePrice[NQ/HQ],
# LevePrice[NQ/HQ], LeveProfit[NQ/HQ]) for the rows whose market-board
# snapshot changed, across all eight crafter sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()

$ws.Range("H112").Value = 2698.125
$ws.Range("J112").Value = 2698.125
$ws.Range("L112").Value = 8094.375
$ws.Range("N112").Value = -10310.375

$ws.Range("H116").Value = 79375.625
$ws.Range("J116").Value = 15000
$ws.Range("L116").Value = 15000
$ws.Range("N116").Value = -21884

$ws.Range("H132").Value = 4349.5
$ws.Range("I132").Value = 3886.1
$ws.Range("J132").Value = 6666.5
$ws.Range("K132").Value = 11658.3
$ws.Range("L132").Value = 19999.5
$ws.Range("M132").Value = -9128.299999999999
$ws.Range("N132").Value = -25059.5

$ws.Range("H137").Value = 2348.9167
$ws.Range("I137").Value = 2493.7
$ws.Range("J137").Value = 1625
$ws.Range("K137").Value = 7481.099999999999
$ws.Range("L137").Value = 4875
$ws.Range("M137").Value = -4931.099999999999
$ws.Range("N137").Value = -9975

$ws.Range("H138").Value = 4717.143
$ws.Range("J138").Value = 4945.9
$ws.Range("L138").Value = 14837.7
$ws.Range("N138").Value = -25117.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6552.625
$ws.Range("I2").Value = 6321.8
$ws.Range("K2").Value = 6321.8
$ws.Range("M2").Value = -6208.8

$ws.Range("H45").Value = 2634.2856
$ws.Range("I45").Value = 2575
$ws.Range("K45").Value = 2575
$ws.Range("M45").Value = -2198

$ws.Range("H61").Value = 4954.2
$ws.Range("I61").Value = 4642.75
$ws.Range("J61").Value = 6200
$ws.Range("K61").Value = 4642.75
$ws.Range("L61").Value = 6200
$ws.Range("M61").Value = -4430.75
$ws.Range("N61").Value = -6624

$ws.Range("H63").Value = 12302.75
$ws.Range("I63").Value = 10127.5
$ws.Range("J63").Value = 14478
$ws.Range("K63").Value = 10127.5
$ws.Range("L63").Value = 14478
$ws.Range("M63").Value = -9441.5
$ws.Range("N63").Value = -15850

$ws.Range("H66").Value = 12302.75
$ws.Range("I66").Value = 10127.5
$ws.Range("J66").Value = 14478
$ws.Range("K66").Value = 50637.5
$ws.Range("L66").Value = 72390
$ws.Range("M66").Value = -47205.5
$ws.Range("N66").Value = -79254

$ws.Range("H74").Value = 4164.6577
$ws.Range("J74").Value = 1098.8334
$ws.Range("L74").Value = 1098.8334
$ws.Range("N74").Value = -2846.8334

$ws.Range("H77").Value = 4164.6577
$ws.Range("J77").Value = 1098.8334
$ws.Range("L77").Value = 5494.166999999999
$ws.Range("N77").Value = -14230.167

$ws.Range("H101").Value = 94944
$ws.Range("J101").Value = 94944
$ws.Range("L101").Value = 94944
$ws.Range("N101").Value = -101434

$ws.Range("H110").Value = 8600
$ws.Range("I110").Value = 3133.3333
$ws.Range("K110").Value = 3133.3333
$ws.Range("M110").Value = -1088.3333

$ws.Range("H116").Value = 6552.625
$ws.Range("I116").Value = 6321.8
$ws.Range("K116").Value = 6321.8
$ws.Range("M116").Value = -4027.8

$ws.Range("H132").Value = 4107.727
$ws.Range("I132").Value = 3568.5
$ws.Range("J132").Value = 9500
$ws.Range("K132").Value = 10705.5
$ws.Range("L132").Value = 28500
$ws.Range("M132").Value = -8175.5
$ws.Range("N132").Value = -33560

$ws.Range("H136").Value = 4954.2
$ws.Range("I136").Value = 4642.75
$ws.Range("J136").Value = 6200
$ws.Range("K136").Value = 13928.25
$ws.Range("L136").Value = 18600
$ws.Range("M136").Value = -11378.25
$ws.Range("N136").Value = -23700

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6552.625
$ws.Range("I3").Value = 6321.8
$ws.Range("K3").Value = 6321.8
$ws.Range("M3").Value = -6207.8

$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 3000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1502
$ws.Range("N99").ClearContents()

$ws.Range("H134").Value = 5044.1113
$ws.Range("I134").Value = 2733.3333
$ws.Range("J134").Value = 9665.666999999999
$ws.Range("K134").Value = 8199.999899999999
$ws.Range("L134").Value = 28997.001
$ws.Range("M134").Value = -5664.999899999999
$ws.Range("N134").Value = -34067.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1349.75
$ws.Range("I31").Value = 1689.6666
$ws.Range("J31").Value = 330
$ws.Range("K31").Value = 1689.6666
$ws.Range("L31").Value = 330
$ws.Range("M31").Value = -1394.6666
$ws.Range("N31").Value = -920

$ws.Range("H34").Value = 1349.75
$ws.Range("I34").Value = 1689.6666
$ws.Range("J34").Value = 330
$ws.Range("K34").Value = 1689.6666
$ws.Range("L34").Value = 330
$ws.Range("M34").Value = -1487.6666
$ws.Range("N34").Value = -734

$ws.Range("H58").Value = 3299.6296
$ws.Range("I58").Value = 3282.6086
$ws.Range("J58").Value = 3397.5
$ws.Range("K58").Value = 3282.6086
$ws.Range("L58").Value = 3397.5
$ws.Range("M58").Value = -3079.6086
$ws.Range("N58").Value = -3803.5

$ws.Range("H99").Value = 2348
$ws.Range("I99").Value = 2348
$ws.Range("K99").Value = 2348
$ws.Range("M99").Value = -850

$ws.Range("H107").Value = 1162.4615
$ws.Range("I107").Value = 870.2
$ws.Range("K107").Value = 870.2
$ws.Range("M107").Value = 1049.8

$ws.Range("H126").Value = 2348
$ws.Range("I126").Value = 2348
$ws.Range("K126").Value = 7044
$ws.Range("M126").Value = -4574

$ws.Range("H134").Value = 3677.5833
$ws.Range("I134").Value = 3921
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 11763
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -9228
$ws.Range("N134").Value = -8070

$ws.Range("H136").Value = 3299.6296
$ws.Range("I136").Value = 3282.6086
$ws.Range("J136").Value = 3397.5
$ws.Range("K136").Value = 9847.825800000001
$ws.Range("L136").Value = 10192.5
$ws.Range("M136").Value = -7297.825800000001
$ws.Range("N136").Value = -15292.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 63.77778
$ws.Range("I2").Value = 50.8
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 304.8
$ws.Range("L2").Value = 480
$ws.Range("M2").Value = -191.8
$ws.Range("N2").Value = -706

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 30000
$ws.Range("J33").Value = 30000
$ws.Range("L33").Value = 30000
$ws.Range("N33").Value = -30504

$ws.Range("H107").Value = 988
$ws.Range("I107").Value = 1140.6
$ws.Range("J107").Value = 733.6667
$ws.Range("K107").Value = 1140.6
$ws.Range("L107").Value = 733.6667
$ws.Range("M107").Value = 779.4000000000001
$ws.Range("N107").Value = -4573.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3625
$ws.Range("I46").Value = 3625
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 3625
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -3437
$ws.Range("N46").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4124.5
$ws.Range("J132").Value = 4124.5
$ws.Range("L132").Value = 12373.5
$ws.Range("N132").Value = -17433.5

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
